$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D10").Value = -8.121599999999994
$ws.Range("D12").Value = -8.092799999999999
$ws.Range("E13").Value = 11.9308
$ws.Range("D18").Value = -8.213299999999995
